$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.546.31'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.83%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.874.04'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.27%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.00%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4734'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.27%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2921'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.67%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06494'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.14%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '22.23'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.78%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07718'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.28%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '97.50'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.45%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7421'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.90%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.870.68'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.00%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.158'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.02%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '274.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.40%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.535.05'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.83%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.40'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.37%  '

$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007520'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.76%  '

$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9991'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.14%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.115.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9996'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.268'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.81%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.184'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.77%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.280'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.38%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.34%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.81'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.82%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.922'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.59%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1001'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.40%  '

$ws.Range("E30").Value = '  -1.68%  '

$ws.Range("E31").Value = '  -0.24%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.295'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.83%  '

$ws.Range("E33").Value = '  +1.80%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04806'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.43%  '

$ws.Range("E35").Value = '  -0.36%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6956'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.29%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.712'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01852'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.48%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.747'
$ws.Range("D39").Style = "Normal"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.230'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.84%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '73.99'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.95%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.969'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.49%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4189'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.30%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9996'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.02%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8336'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.98%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '102.67'
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.295'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.38%  '

$ws.Range("E48").Value = '  +1.18%  '

$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '930.34'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.78%  '

$ws.Range("B50").Value = 'Aptos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.985'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.34%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05647'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.36%  '
